$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$b64_2 = "W0NyYWlnJUJyYW1tZXIlTlVMTCUxLCAgICAgICAgICAgICAgICAgU3RhY3klRG9ub2h1ZSVOVUxMJTEsICAgICAgICAgICAgICAgICBUaW1vdGh5JUVsd2VsbCVOVUxMJTEsICAgICAgICAgICAgICAgICBFbGl6YSVGaXNoYmVpbiVOVUxMJTEsICAgICAgICAgICAgICAgICBEJ0FudGUlRm9yc2NoaW5vJU5VTEwlMSwgICAgICAgICAgICAgICAgIERvcm90aHklSG9ybmUlTlVMTCUxLCAgICAgICAgICAgICAgICAgQnVmZnklTGxveWQtS3JlamNpJU5VTEwlMSwgICAgICAgICAgICAgICAgIEplc3NpY2ElTGl0dGxlJU5VTEwlMSwgICAgICAgICAgICAgICAgIEJpc3RyYSVOaWtpZm9yb3ZhJU5VTEwlMSwgICAgICAgICAgICAgICAgIEVsaXphYmV0aCVXaW50ZXJiYXVlciVOVUxMJTFd"
$bytes_2 = [System.Convert]::FromBase64String($b64_2)
$str_2 = [System.Text.Encoding]::UTF8.GetString($bytes_2)
$ws.Range("E2").Value = $str_2

$b64_3 = "W0F5ZXNoYSVBcHBhJUF5ZXNoYS5hcHBhQHVjc2YuZWR1JTEsICAgICAgICAgICAgICAgICBHYWJyaWVsJUNoYW1pZSVOVUxMJTEsICAgICAgICAgICAgICAgICBBZW5vciVTYXd5ZXIlTlVMTCUxLCAgICAgICAgICAgICAgICAgS2ltYmVybHklQmFsdHplbGwlTlVMTCUxLCAgICAgICAgICAgICAgICAgS2F0aHJ5biVEaXBwZWxsJU5VTEwlMSwgICAgICAgICAgICAgICAgIFNhbHUlUmliZWlybyVOVUxMJTEsICAgICAgICAgICAgICAgICBFbGlhcyVEdWFydGUlTlVMTCUxLCAgICAgICAgICAgICAgICAgSm9hbm5hJVZpbmRlbiVOVUxMJTEsICAgICAgICAgICAgICAgICBDTElBSFVCJUNvbnNvcnRpdW0lTlVMTCUxLCAgICAgICAgICAgICAgICAgSm9uYXRoYW4lS3JhbWVyLUZlbGRtYW4lTlVMTCUxLCAgICAgICAgICAgICAgICAgU2hhaHJ5YXIlUmFoZGFyaSVOVUxMJTEsICAgICAgICAgICAgICAgICBEb3VnJU1hY0ludG9zaCVOVUxMJTEsICAgICAgICAgICAgICAgICBLYXRoZXJpbmUlTmljaG9sc29uJU5VTEwlMSwgICAgICAgICAgICAgICAgIEpvbmF0aGFuJUltJU5VTEwlMSwgICAgICAgICAgICAgICAgIERpYW5lJUhhdmxpciVOVUxMJTEsICAgICAgICAgICAgICAgICBCcnlhbiVHcmVlbmhvdXNlJU5VTEwlMV0="
$bytes_3 = [System.Convert]::FromBase64String($b64_3)
$str_3 = [System.Text.Encoding]::UTF8.GetString($bytes_3)
$ws.Range("E3").Value = $str_3

$b64_4 = "W01pbiBDaGVvbCVDaGFuZyVOVUxMJTEsICAgICAgICAgICAgICAgICBXYW4tU2VvayVTZW8lTlVMTCUxLCAgICAgICAgICAgICAgICAgRG9uZ2h3aSVQYXJrJU5VTEwlMSwgICAgICAgICAgICAgICAgIEppYW4lSHVyJU5VTEwlMF0="
$bytes_4 = [System.Convert]::FromBase64String($b64_4)
$str_4 = [System.Text.Encoding]::UTF8.GetString($bytes_4)
$ws.Range("E4").Value = $str_4

$b64_5 = "W01hcmNpIEwuJURyZWVzJU5VTEwlMSwgICAgICAgICAgICAgICAgIE1pYSBBLiVQYXBhcyVOVUxMJTEsICAgICAgICAgICAgICAgICBUZXJyaSBFLiVDb3JibyVOVUxMJTEsICAgICAgICAgICAgICAgICBLaW1iZXJseSBELiVXaWxsaWFtcyVOVUxMJTEsICAgICAgICAgICAgICAgICBTaGFyb24gVC4lS3VyZnVlcnN0JU5VTEwlMiwgICAgICAgICAgICAgICAgIFNoYXJvbiBULiVLdXJmdWVyc3QlTlVMTCUwXQ=="
$bytes_5 = [System.Convert]::FromBase64String($b64_5)
$str_5 = [System.Text.Encoding]::UTF8.GetString($bytes_5)
$ws.Range("E5").Value = $str_5

$b64_6 = "W0VyaW4gRi4lRmx5bm4lZmx5bm5lQGVtYWlsLmNob3AuZWR1JTEsICAgICAgICAgICAgICAgICBFbGl6YWJldGglS3VobiVOVUxMJTEsICAgICAgICAgICAgICAgICBNb2hhbW1lZCVTaGFpayVOVUxMJTEsICAgICAgICAgICAgICAgICBFbGl6YWJldGglVGFyciVOVUxMJTEsICAgICAgICAgICAgICAgICBOaWNvbGUlU2NhdHRvbGluaSVOVUxMJTEsICAgICAgICAgICAgICAgICBBbGxpc29uJUJhbGxhbnRpbmUlTlVMTCUxXQ=="
$bytes_6 = [System.Convert]::FromBase64String($b64_6)
$str_6 = [System.Text.Encoding]::UTF8.GetString($bytes_6)
$ws.Range("E6").Value = $str_6

$b64_7 = "W1Njb3R0IEEuJUdvbGRiZXJnJU5VTEwlMSwgICAgICAgICAgICAgICAgIFJvYmVydCBBLiVCb25hY2NpJU5VTEwlMSwgICAgICAgICAgICAgICAgIEx1Y2FzIEMuJUNhcmxzb24lTlVMTCUxLCAgICAgICAgICAgICAgICAgQ2hhcmxlcyBULiVQdSVOVUxMJTEsICAgICAgICAgICAgICAgICBDaHJpc3RpbmUgUy4lUml0Y2hpZSVOVUxMJTFd"
$bytes_7 = [System.Convert]::FromBase64String($b64_7)
$str_7 = [System.Text.Encoding]::UTF8.GetString($bytes_7)
$ws.Range("E7").Value = $str_7

$b64_8 = "W1RyYXZpcyVTYW5jaGV6JU5VTEwlMCwgICAgICAgICAgICAgICAgIFNhZGh1JVBhbmRhJU5VTEwlMiwgICAgICAgICAgICAgICAgIFNhZGh1JVBhbmRhJU5VTEwlMCwgICAgICAgICAgICAgICAgIEVicmFoaW0lS2hhamVoJU5VTEwlMSwgICAgICAgICAgICAgICAgIEFsZXhhbmRyYSVIYWxhbGF1JWFsZXhhbmRyYS5oYWxhbGF1QGJlYXVtb250LmVkdSUyLCAgICAgICAgICAgICAgICAgQWxleGFuZHJhJUhhbGFsYXUlYWxleGFuZHJhLmhhbGFsYXVAYmVhdW1vbnQuZWR1JTAsICAgICAgICAgICAgICAgICBKZWZmcmV5JURpdGtvZmYlTlVMTCUyLCAgICAgICAgICAgICAgICAgSmVmZnJleSVEaXRrb2ZmJU5VTEwlMCwgICAgICAgICAgICAgICAgIEplc3NpY2ElSGFtaWx0b24lTlVMTCUyLCAgICAgICAgICAgICAgICAgSmVzc2ljYSVIYW1pbHRvbiVOVUxMJTAsICAgICAgICAgICAgICAgICBBcnlhbmElU2hhcnJhayVOVUxMJTIsICAgICAgICAgICAgICAgICBBcnlhbmElU2hhcnJhayVOVUxMJTAsICAgICAgICAgICAgICAgICBBaW1lbiVWYW5vb2QlTlVMTCUyLCAgICAgICAgICAgICAgICAgQWltZW4lVmFub29kJU5VTEwlMCwgICAgICAgICAgICAgICAgIEFtciVBYmJhcyVOVUxMJTIsICAgICAgICAgICAgICAgICBBbXIlQWJiYXMlTlVMTCUwLCAgICAgICAgICAgICAgICAgSmFtZXMlWmlhZGVoJU5VTEwlMiwgICAgICAgICAgICAgICAgIEphbWVzJVppYWRlaCVOVUxMJTBd"
$bytes_8 = [System.Convert]::FromBase64String($b64_8)
$str_8 = [System.Text.Encoding]::UTF8.GetString($bytes_8)
$ws.Range("E8").Value = $str_8

$b64_10 = "W1NhbmcgSWwlS2ltJU5VTEwlMSwgICAgICAgICAgICAgICAgIEppIFlvbmclTGVlJU5VTEwlMiwgICAgICAgICAgICAgICAgIEppIFlvbmclTGVlJU5VTEwlMF0="
$bytes_10 = [System.Convert]::FromBase64String($b64_10)
$str_10 = [System.Text.Encoding]::UTF8.GetString($bytes_10)
$ws.Range("E10").Value = $str_10

$b64_11 = "W0thdGhlcmluZSVIaWxsJXhyZWYgbm8gZW1haWwlMSwgICAgICAgICAgUm9ieW4lQ2FtcGJlbGwleHJlZiBubyBlbWFpbCUxLCAgICAgICAgICBDYWxsdW0lTXV0Y2gleHJlZiBubyBlbWFpbCUxLCAgICAgICAgICBPbGl2ZXIlS29jaCV4cmVmIG5vIGVtYWlsJTEsICAgICAgICAgIENsYWlyZSVNYWNraW50b3NoJXhyZWYgbm8gZW1haWwlMV0="
$bytes_11 = [System.Convert]::FromBase64String($b64_11)
$str_11 = [System.Text.Encoding]::UTF8.GetString($bytes_11)
$ws.Range("E11").Value = $str_11

$b64_12 = "W0tpIFRhZSVLd29uJU5VTEwlMSwgICAgICAgICAgICAgICAgIEphZS1Ib29uJUtvJU5VTEwlMiwgICAgICAgICAgICAgICAgIEphZS1Ib29uJUtvJU5VTEwlMCwgICAgICAgICAgICAgICAgIEhlZWp1biVTaGluJU5VTEwlMiwgICAgICAgICAgICAgICAgIEhlZWp1biVTaGluJU5VTEwlMCwgICAgICAgICAgICAgICAgIE1pbmtpJVN1bmclTlVMTCUyLCAgICAgICAgICAgICAgICAgTWlua2klU3VuZyVOVUxMJTAsICAgICAgICAgICAgICAgICBKaW4gWW9uZyVLaW0lTlVMTCU0LCAgICAgICAgICAgICAgICAgSmluIFlvbmclS2ltJU5VTEwlMF0="
$bytes_12 = [System.Convert]::FromBase64String($b64_12)
$str_12 = [System.Text.Encoding]::UTF8.GetString($bytes_12)
$ws.Range("E12").Value = $str_12

$b64_13 = "W0VsbHklTGVlJU5VTEwlMSwgICAgICAgICAgICAgICAgIE51cnVsIFlhcWVlbiVNb2hkIEVzYSVOVUxMJTEsICAgICAgICAgICAgICAgICBUb25nIE1pbmclV2VlJU5VTEwlMSwgICAgICAgICAgICAgICAgIENodW4gSWFuJVNvbyVOVUxMJTFd"
$bytes_13 = [System.Convert]::FromBase64String($b64_13)
$str_13 = [System.Text.Encoding]::UTF8.GetString($bytes_13)
$ws.Range("E13").Value = $str_13

$b64_14 = "W1lvbmcgU2hpayVLd29uJU5VTEwlMSwgICAgICAgICAgICAgICAgIFN1biBIeW8lUGFyayVOVUxMJTIsICAgICAgICAgICAgICAgICBTdW4gSHlvJVBhcmslTlVMTCUwLCAgICAgICAgICAgICAgICAgSHl1biBKdW5nJUtpbSVOVUxMJTIsICAgICAgICAgICAgICAgICBIeXVuIEp1bmclS2ltJU5VTEwlMCwgICAgICAgICAgICAgICAgIEppIFllb24lTGVlJU5VTEwlMCwgICAgICAgICAgICAgICAgIEppIFllb24lTGVlJU5VTEwlMCwgICAgICAgICAgICAgICAgIE1pLXJpJUh5dW4lTlVMTCUyLCAgICAgICAgICAgICAgICAgTWktcmklSHl1biVOVUxMJTAsICAgICAgICAgICAgICAgICBIeXVuIGFoJUtpbSVOVUxMJTIsICAgICAgICAgICAgICAgICBIeXVuIGFoJUtpbSVOVUxMJTAsICAgICAgICAgICAgICAgICBKYWUgU2VvayVQYXJrJU5VTEwlMCwgICAgICAgICAgICAgICAgIEphZSBTZW9rJVBhcmslTlVMTCUwXQ=="
$bytes_14 = [System.Convert]::FromBase64String($b64_14)
$str_14 = [System.Text.Encoding]::UTF8.GetString($bytes_14)
$ws.Range("E14").Value = $str_14

$b64_15 = "W1BvLVRpbmclTGluJU5VTEwlMSwgICAgICAgICAgICAgICAgIFRpbmctWXVhbiVOaSVOVUxMJTEsICAgICAgICAgICAgICAgICBUcmVuLVlpJUNoZW4lTlVMTCUxLCAgICAgICAgICAgICAgICAgQ2hpaC1QZWklU3UlTlVMTCUxLCAgICAgICAgICAgICAgICAgSHNpYW8tRmVuJVN1biVOVUxMJTEsICAgICAgICAgICAgICAgICBNdS1LdWFuJUNoZW4lTlVMTCUxLCAgICAgICAgICAgICAgICAgQ2h1LUNodW5nJUNob3UlTlVMTCUxLCAgICAgICAgICAgICAgICAgUG8tWXUlV2FuZyVOVUxMJTEsICAgICAgICAgICAgICAgICBZYW4tUmVuJUxpbiVINjIxMy5sYWNAZ21haWwuY29tJTFd"
$bytes_15 = [System.Convert]::FromBase64String($b64_15)
$str_15 = [System.Text.Encoding]::UTF8.GetString($bytes_15)
$ws.Range("E15").Value = $str_15

$b64_16 = "W0RhdmlkIEElTGluZGhvbG0lZGF2aWQuYS5saW5kaG9sbTQubWlsQG1haWwubWlsJTEsICAgICAgICAgICAgICAgICBKb2huIEwlS2lsZXklTlVMTCUyLCAgICAgICAgICAgICAgICAgSm9obiBMJUtpbGV5JU5VTEwlMCwgICAgICAgICAgICAgICAgIE5hdGhhbiBLJUphbnNlbiVOVUxMJTEsICAgICAgICAgICAgICAgICBSb2JlcnQgVCVIb2FyZCVOVUxMJTEsICAgICAgICAgICAgICAgICBNYXR0aGV3IFIlQm9uZGFyeWslTlVMTCUxLCAgICAgICAgICAgICAgICAgRWxpemFiZXRoIE0lU3RhbmxleSVOVUxMJTEsICAgICAgICAgICAgICAgICBHYWRpZWwgUiVBbHZhcmFkbyVOVUxMJTEsICAgICAgICAgICAgICAgICBBbmEgRSVNYXJrZWx6JU5VTEwlMSwgICAgICAgICAgICAgICAgIFJvYmVydCBKJUN5YnVsc2tpJU5VTEwlMSwgICAgICAgICAgICAgICAgIEphc29uIEYlT2t1bGljeiVOVUxMJTFd"
$bytes_16 = [System.Convert]::FromBase64String($b64_16)
$str_16 = [System.Text.Encoding]::UTF8.GetString($bytes_16)
$ws.Range("E16").Value = $str_16

$b64_17 = "W0NoYXJtYWluZSBNYWxlbmFiJU1hbmF1aXMlTlVMTCUxLCAgICAgICAgICAgICAgICAgTWFydmluJUxvaCVtYXJ2aW4ubG9oQG1vaGguY29tLnNnJTEsICAgICAgICAgICAgICAgICBKYW1lcyVLd2FuJU5VTEwlMSwgICAgICAgICAgICAgICAgIEpvaG4lQ2h1YSBNaW5nemhvdSVOVUxMJTEsICAgICAgICAgICAgICAgICBIYW4gSmllJVRlbyVOVUxMJTEsICAgICAgICAgICAgICAgICBEYXZpZCVUZW5nIEt1YW4gUGVuZyVOVUxMJTEsICAgICAgICAgICAgICAgICBTaGF3biVWYXNvbyBTdXNoaWxhbiVOVUxMJTEsICAgICAgICAgICAgICAgICBZZWUgU2luJUxlbyVOVUxMJTEsICAgICAgICAgICAgICAgICBBbmclSG91JU5VTEwlMV0="
$bytes_17 = [System.Convert]::FromBase64String($b64_17)
$str_17 = [System.Text.Encoding]::UTF8.GetString($bytes_17)
$ws.Range("E17").Value = $str_17

$b64_18 = "W1JpdmtlZXMlU2NvdHQgQS4lY29yZUdpdmVzTm9FbWFpbCUxLCAgICAgICAgICAgICAgICBSb2JlcnNvbiVTaGFtYXJpYWwlY29yZUdpdmVzTm9FbWFpbCUxXQ=="
$bytes_18 = [System.Convert]::FromBase64String($b64_18)
$str_18 = [System.Text.Encoding]::UTF8.GetString($bytes_18)
$ws.Range("E18").Value = $str_18

$b64_19 = "W0V1bmh5ZSVTZW8lTlVMTCUxLCAgICAgICAgICAgICAgICAgRXVuY2hhbiVNdW4lTlVMTCUyLCAgICAgICAgICAgICAgICAgRXVuY2hhbiVNdW4lTlVMTCUwLCAgICAgICAgICAgICAgICAgV29uc29vbCVLaW0lTlVMTCUyLCAgICAgICAgICAgICAgICAgV29uc29vbCVLaW0lTlVMTCUwLCAgICAgICAgICAgICAgICAgQ2hhbmdod2FuJUxlZSVOVUxMJTIsICAgICAgICAgICAgICAgICBDaGFuZ2h3YW4lTGVlJU5VTEwlMF0="
$bytes_19 = [System.Convert]::FromBase64String($b64_19)
$str_19 = [System.Text.Encoding]::UTF8.GetString($bytes_19)
$ws.Range("E19").Value = $str_19

$b64_20 = "W0FkaXR5YSVTaGFoJXNoYWguYWRpdHlhQG1heW8uZWR1JTEsICAgICAgICAgICAgICAgICBEb3VnbGFzJUNoYWxsZW5lciVOVUxMJTEsICAgICAgICAgICAgICAgICBBYXJvbiBKLiVUYW5kZSVOVUxMJTEsICAgICAgICAgICAgICAgICBNYXJ5YW0lTWFobW9vZCVOVUxMJTEsICAgICAgICAgICAgICAgICBKb2huIEMuJU/igJlIb3JvJU5VTEwlMSwgICAgICAgICAgICAgICAgIEVsaWUlQmVyYmFyaSVOVUxMJTEsICAgICAgICAgICAgICAgICBTYXJhaCBKLiVDcmFuZSVOVUxMJTFd"
$bytes_20 = [System.Convert]::FromBase64String($b64_20)
$str_20 = [System.Text.Encoding]::UTF8.GetString($bytes_20)
$ws.Range("E20").Value = $str_20

$b64_21 = "W0FuZ2llIE4uJVRvbiVOVUxMJTEsICAgICAgICAgICAgICAgICBUYXJhbmclSmV0aHdhJU5VTEwlMSwgICAgICAgICAgICAgICAgIEthcmVuJVdhdGVycyVOVUxMJTIsICAgICAgICAgICAgICAgICBLYXJlbiVXYXRlcnMlTlVMTCUwLCAgICAgICAgICAgICAgICAgTGVpZ2ggTC4lU3BlaWNoZXIlTlVMTCUyLCAgICAgICAgICAgICAgICAgTGVpZ2ggTC4lU3BlaWNoZXIlTlVMTCUwLCAgICAgICAgICAgICAgICAgRGF3biVGcmFuY2lzJU5VTEwlMiwgICAgICAgICAgICAgICAgIERhd24lRnJhbmNpcyVOVUxMJTBd"
$bytes_21 = [System.Convert]::FromBase64String($b64_21)
$str_21 = [System.Text.Encoding]::UTF8.GetString($bytes_21)
$ws.Range("E21").Value = $str_21

$b64_22 = "W0thdGUlTWFyayVOVUxMJTEsICAgICAgICAgICAgICAgICBLYXRpZSVTdGVlbCVOVUxMJTEsICAgICAgICAgICAgICAgICBKYW5ldCVTdGV2ZW5zb24lTlVMTCUxLCAgICAgICAgICAgICAgICAgQ2hyaXN0aW5lJUV2YW5zJU5VTEwlMSwgICAgICAgICAgICAgICAgIER1bmNhbiVNY0Nvcm1pY2slTlVMTCUxLCAgICAgICAgICAgICAgICAgTG9ybmElV2lsbG9ja3MlTlVMTCUxLCAgICAgICAgICAgICAgICAgQWxpc29uJU1jQ2FsbHVtJU5VTEwlMSwgICAgICAgICAgICAgICAgIExhdXJhJUpvbmVzJU5VTEwlMSwgICAgICAgICAgICAgICAgIEluZ29sZnVyJUpvaGFubmVzc2VuJU5VTEwlMSwgICAgICAgICAgICAgICAgIEthdGUlVGVtcGxldG9uJU5VTEwlMSwgICAgICAgICAgICAgICAgIE9saXZlciVLb2NoJU5VTEwlMCwgICAgICAgICAgICAgICAgIENsYWlyZSVNYWNraW50b3NoJU5VTEwlMV0="
$bytes_22 = [System.Convert]::FromBase64String($b64_22)
$str_22 = [System.Text.Encoding]::UTF8.GetString($bytes_22)
$ws.Range("E22").Value = $str_22
